$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 21 - new entry: "46. Permutations" (Backtracking category, Medium difficulty)
$ws.Range("A21").Value = "Backtracking"

$ws.Range("B21").Value = "Medium"
$ws.Range("B21").Style = "Neutral"

$ws.Range("C21").Value = "46. Permutations"
$ws.Hyperlinks.Add($ws.Range("C21"), "https://leetcode.com/problems/permutations/")
$ws.Range("C21").Style = "Neutral"

$noteText = "Break down (Recursive ans) or Build up (Iterative ans) to the answer. Permutations for [1] is just [[1]] right, for [1,2] its [[1,2],[2,1]], which is just permutation of [1] but with 2 inserted at all possible indices, and same for [1,2,3] and so on.`n(Also do not modify any element of res directly, make a copy first)`nRecursive - have base case be nums == 0: return [[]]. Then recusively get permutations for nums[1:] and store it in perms, after that just insert nums[0] at all possible positions, store it in a res variable and return res.`nIterative - start with res = [[]]. Outermost loop loops over all nums, then go over each list in res, and then go over each index and insert the selected num at all possible positions. Append to a temp array since we cannot modify res while looping over it and then replace res with temp and return res in the end."
$ws.Range("D21").Value = $noteText
$d21 = $ws.Range("D21")
$d21.Characters(320,9).Font.Bold = $true
$d21.Characters(330,1).Font.Bold = $true
$d21.Characters(544,12).Font.Bold = $true

$ws.Rows.Item(21).RowHeight = 115.2

$ws.Range("D21").Select()
